$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2490446666666667
$ws.Range("H2").Value = 0.747134
$ws.Range("I2").Value = 0.00222126824250593
$ws.Range("J2").Value = 0.00222126824250593
$ws.Range("M2").Value = 28.19948866666667
$ws.Range("N2").Value = 84.598466
$ws.Range("O2").Value = 0.7357427920402423
$ws.Range("P2").Value = 0.7357427920402422
$ws.Range("Q2").Value = 7.022932255160445
$ws.Range("R2").Value = 63.206390296444
$ws.Range("S2").Value = 0.001634282098611635
$ws.Range("T2").Value = 0.001634282098611635

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2490446666666667
$ws.Range("H3").Value = 0.747134
$ws.Range("I3").Value = 0.00222126824250593
$ws.Range("J3").Value = 0.00222126824250593
$ws.Range("O3").Value = 0.2029336910395279
$ws.Range("P3").Value = 0.2029336910395278
$ws.Range("Q3").Value = 1.937075809479778
$ws.Range("R3").Value = 17.433682285318
$ws.Range("S3").Value = 0.0004507701632406134
$ws.Range("T3").Value = 0.0004507701632406134

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2490446666666667
$ws.Range("H4").Value = 0.747134
$ws.Range("I4").Value = 0.00222126824250593
$ws.Range("J4").Value = 0.00222126824250593
$ws.Range("M4").Value = 2.350402666666667
$ws.Range("N4").Value = 7.051208000000001
$ws.Range("O4").Value = 0.0613235169202299
$ws.Range("P4").Value = 0.06132351692022989
$ws.Range("Q4").Value = 0.5853552486524446
$ws.Range("R4").Value = 5.268197237872
$ws.Range("S4").Value = 0.0001362159806536817
$ws.Range("T4").Value = 0.0001362159806536817

$ws.Range("I5").Value = 0.9377915177839022
$ws.Range("J5").Value = 0.9377915177839021
$ws.Range("M5").Value = 28.19948866666667
$ws.Range("N5").Value = 84.598466
$ws.Range("O5").Value = 0.7357427920402423
$ws.Range("P5").Value = 0.7357427920402422
$ws.Range("Q5").Value = 2964.993679210202
$ws.Range("R5").Value = 26684.94311289181
$ws.Range("S5").Value = 0.6899733496459848
$ws.Range("T5").Value = 0.6899733496459846

$ws.Range("I6").Value = 0.9377915177839022
$ws.Range("J6").Value = 0.9377915177839021
$ws.Range("O6").Value = 0.2029336910395279
$ws.Range("P6").Value = 0.2029336910395278
$ws.Range("S6").Value = 0.1903094941294483
$ws.Range("T6").Value = 0.1903094941294483

$ws.Range("I7").Value = 0.9377915177839022
$ws.Range("J7").Value = 0.9377915177839021
$ws.Range("M7").Value = 2.350402666666667
$ws.Range("N7").Value = 7.051208000000001
$ws.Range("O7").Value = 0.0613235169202299
$ws.Range("P7").Value = 0.06132351692022989
$ws.Range("Q7").Value = 247.1296246766036
$ws.Range("R7").Value = 2224.166622089433
$ws.Range("S7").Value = 0.05750867400846921
$ws.Range("T7").Value = 0.05750867400846919

$ws.Range("G8").Value = 6.725660333333334
$ws.Range("I8").Value = 0.05998721397359182
$ws.Range("J8").Value = 0.05998721397359182
$ws.Range("M8").Value = 28.19948866666667
$ws.Range("N8").Value = 84.598466
$ws.Range("O8").Value = 0.7357427920402423
$ws.Range("P8").Value = 0.7357427920402422
$ws.Range("Q8").Value = 189.6601823456829
$ws.Range("R8").Value = 1706.941641111146
$ws.Range("S8").Value = 0.04413516029564588
$ws.Range("T8").Value = 0.04413516029564588

$ws.Range("G9").Value = 6.725660333333334
$ws.Range("I9").Value = 0.05998721397359182
$ws.Range("J9").Value = 0.05998721397359182
$ws.Range("O9").Value = 0.2029336910395279
$ws.Range("P9").Value = 0.2029336910395278
$ws.Range("Q9").Value = 52.31235869794856
$ws.Range("R9").Value = 470.8112282815371
$ws.Range("S9").Value = 0.01217342674683893
$ws.Range("T9").Value = 0.01217342674683893

$ws.Range("G10").Value = 6.725660333333334
$ws.Range("I10").Value = 0.05998721397359182
$ws.Range("J10").Value = 0.05998721397359182
$ws.Range("M10").Value = 2.350402666666667
$ws.Range("N10").Value = 7.051208000000001
$ws.Range("O10").Value = 0.0613235169202299
$ws.Range("P10").Value = 0.06132351692022989
$ws.Range("Q10").Value = 15.80800998256089
$ws.Range("R10").Value = 142.272089843048
$ws.Range("S10").Value = 0.003678626931107009
$ws.Range("T10").Value = 0.003678626931107008
